# Add Trade #4 (closed OPEN state, leadlag DOWN) to the "All Trades" and
# "leadlag" worksheets at row 5.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "leadlag")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Cells.Item(5, 1).Value = 4
    $ws.Cells.Item(5, 2).Value = "2026-02-16"
    $ws.Cells.Item(5, 3).Value = "21:51:15"
    $ws.Cells.Item(5, 4).Value = "leadlag"
    $ws.Cells.Item(5, 5).Value = "DOWN"
    $ws.Cells.Item(5, 6).Value = 68298.795
    $ws.Cells.Item(5, 7).Value = ""
    $ws.Cells.Item(5, 8).Value = "OPEN"
    $ws.Cells.Item(5, 9).Value = 0
    $ws.Cells.Item(5, 10).Value = 0
    $ws.Cells.Item(5, 11).Value = 100
    $ws.Cells.Item(5, 12).Value = 0.75
    $ws.Cells.Item(5, 13).Value = "Coinbase leading with -0.143% move"
    $ws.Cells.Item(5, 14).Value = ""
    $ws.Cells.Item(5, 15).Value = 0
}
